$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 8.5
$ws.Range("J3").Value = 8.5
$ws.Range("K3").Value = 2.4
$ws.Range("N3").Value = 12
$ws.Range("Q3").Value = 1.75
$ws.Range("R3").Value = 2.05
$ws.Range("U3").Value = 2.2
$ws.Range("V3").Value = 1.62

# Row 6 updates
$ws.Range("N6").Value = 6.75
